# Update "BVTStL" sheet: mark aircraft, rail, and ships (rows 4-6) as
# subject to LCFS (Boolean 0 -> 1) in both the LDV-fuel and HDV-fuel
# columns (B and C), and move the sheet's saved selection to D7.

$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("BVTStL")

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1

# Move the cursor on the BVTStL sheet to D7 (this briefly activates the
# sheet, same as it would in the real Excel UI).
$ws.Range("D7").Select()

# Restore "About" as the active/visible tab, matching the original
# workbook view state (only BVTStL's stored selection should change).
$about.Activate()
